$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("Lesson 4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Start = $r.End - 1
$r.Bold = 1
$r.Text = "7"
$r.Bold = 0
Write-Output "WordOpenXML:"
Write-Output $r.WordOpenXML
